$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "Symptom – Respiratory"
$ws.Cells.Item(2, 3).Value = 1970
$ws.Cells.Item(2, 5).Value = 27.1
$ws.Cells.Item(2, 6).Value = 26.08
$ws.Cells.Item(2, 7).Value = 28.12

$ws.Cells.Item(3, 2).Value = "Injuries & adverse effects"
$ws.Cells.Item(3, 3).Value = 1155
$ws.Cells.Item(3, 5).Value = 15.89
$ws.Cells.Item(3, 6).Value = 15.05
$ws.Cells.Item(3, 7).Value = 16.73

$ws.Cells.Item(4, 2).Value = "Symptom – Digestive"
$ws.Cells.Item(4, 3).Value = 878
$ws.Cells.Item(4, 5).Value = 12.08
$ws.Cells.Item(4, 6).Value = 11.33
$ws.Cells.Item(4, 7).Value = 12.83

$ws.Cells.Item(5, 2).Value = "Symptom – Nervous"
$ws.Cells.Item(5, 3).Value = 737
$ws.Cells.Item(5, 5).Value = 10.14
$ws.Cells.Item(5, 6).Value = 9.44
$ws.Cells.Item(5, 7).Value = 10.83

$ws.Cells.Item(6, 2).Value = "Symptom – Circulatory"
$ws.Cells.Item(6, 3).Value = 718
$ws.Cells.Item(6, 5).Value = 9.880000000000001
$ws.Cells.Item(6, 6).Value = 9.19
$ws.Cells.Item(6, 7).Value = 10.56

$ws.Cells.Item(7, 2).Value = "Other"
$ws.Cells.Item(7, 3).Value = 464
$ws.Cells.Item(7, 5).Value = 6.38
$ws.Cells.Item(7, 6).Value = 5.82
$ws.Cells.Item(7, 7).Value = 6.94

$ws.Cells.Item(8, 2).Value = "Diseases (patient-stated)"
$ws.Cells.Item(8, 3).Value = 453
$ws.Cells.Item(8, 5).Value = 6.23
$ws.Cells.Item(8, 6).Value = 5.68
$ws.Cells.Item(8, 7).Value = 6.79

$ws.Cells.Item(9, 2).Value = "Symptom – General"
$ws.Cells.Item(9, 3).Value = 302
$ws.Cells.Item(9, 5).Value = 4.15
$ws.Cells.Item(9, 6).Value = 3.7
$ws.Cells.Item(9, 7).Value = 4.61

$ws.Cells.Item(10, 2).Value = "Symptom – Skin/Hair/Nails"
$ws.Cells.Item(10, 3).Value = 211
$ws.Cells.Item(10, 5).Value = 2.9
$ws.Cells.Item(10, 6).Value = 2.52
$ws.Cells.Item(10, 7).Value = 3.29

$ws.Cells.Item(11, 2).Value = "Uncodable/Unknown"
$ws.Cells.Item(11, 3).Value = 192
$ws.Cells.Item(11, 5).Value = 2.64
$ws.Cells.Item(11, 6).Value = 2.27
$ws.Cells.Item(11, 7).Value = 3.01

$ws.Cells.Item(12, 2).Value = "Symptom – Musculoskeletal"
$ws.Cells.Item(12, 3).Value = 190
$ws.Cells.Item(12, 5).Value = 2.61
$ws.Cells.Item(12, 6).Value = 2.25
$ws.Cells.Item(12, 7).Value = 2.98

$ws.Cells.Item(13, 2).Value = "Symptom – Respiratory"
$ws.Cells.Item(13, 3).Value = 988
$ws.Cells.Item(13, 5).Value = 49.82
$ws.Cells.Item(13, 6).Value = 47.62
$ws.Cells.Item(13, 7).Value = 52.02

$ws.Cells.Item(14, 2).Value = "Symptom – Nervous"
$ws.Cells.Item(14, 3).Value = 223
$ws.Cells.Item(14, 5).Value = 11.25
$ws.Cells.Item(14, 6).Value = 9.859999999999999
$ws.Cells.Item(14, 7).Value = 12.64

$ws.Cells.Item(15, 2).Value = "Injuries & adverse effects"
$ws.Cells.Item(15, 3).Value = 158
$ws.Cells.Item(15, 5).Value = 7.97
$ws.Cells.Item(15, 6).Value = 6.78
$ws.Cells.Item(15, 7).Value = 9.16

$ws.Cells.Item(16, 2).Value = "Symptom – Digestive"
$ws.Cells.Item(16, 3).Value = 143
$ws.Cells.Item(16, 5).Value = 7.21
$ws.Cells.Item(16, 6).Value = 6.07
$ws.Cells.Item(16, 7).Value = 8.35

$ws.Cells.Item(17, 2).Value = "Symptom – Circulatory"
$ws.Cells.Item(17, 3).Value = 133
$ws.Cells.Item(17, 5).Value = 6.71
$ws.Cells.Item(17, 6).Value = 5.61
$ws.Cells.Item(17, 7).Value = 7.81

$ws.Cells.Item(18, 2).Value = "Other"
$ws.Cells.Item(18, 3).Value = 113
$ws.Cells.Item(18, 5).Value = 5.7
$ws.Cells.Item(18, 6).Value = 4.68
$ws.Cells.Item(18, 7).Value = 6.72

$ws.Cells.Item(19, 2).Value = "Symptom – General"
$ws.Cells.Item(19, 3).Value = 82
$ws.Cells.Item(19, 5).Value = 4.14
$ws.Cells.Item(19, 6).Value = 3.26
$ws.Cells.Item(19, 7).Value = 5.01

$ws.Cells.Item(20, 2).Value = "Diseases (patient-stated)"
$ws.Cells.Item(20, 3).Value = 60
$ws.Cells.Item(20, 5).Value = 3.03
$ws.Cells.Item(20, 6).Value = 2.27
$ws.Cells.Item(20, 7).Value = 3.78

$ws.Cells.Item(21, 2).Value = "Symptom – Skin/Hair/Nails"
$ws.Cells.Item(21, 3).Value = 33
$ws.Cells.Item(21, 5).Value = 1.66
$ws.Cells.Item(21, 6).Value = 1.1
$ws.Cells.Item(21, 7).Value = 2.23

$ws.Cells.Item(22, 2).Value = "Symptom – Musculoskeletal"
$ws.Cells.Item(22, 3).Value = 28
$ws.Cells.Item(22, 5).Value = 1.41
$ws.Cells.Item(22, 6).Value = 0.89
$ws.Cells.Item(22, 7).Value = 1.93

$ws.Cells.Item(23, 2).Value = "Uncodable/Unknown"
$ws.Cells.Item(23, 3).Value = 22
$ws.Cells.Item(23, 5).Value = 1.11
$ws.Cells.Item(23, 6).Value = 0.65
$ws.Cells.Item(23, 7).Value = 1.57

$ws.Cells.Item(24, 2).Value = "Symptom – Respiratory"
$ws.Cells.Item(24, 3).Value = 460
$ws.Cells.Item(24, 5).Value = 34.07
$ws.Cells.Item(24, 6).Value = 31.55
$ws.Cells.Item(24, 7).Value = 36.6

$ws.Cells.Item(25, 2).Value = "Symptom – Digestive"
$ws.Cells.Item(25, 3).Value = 163
$ws.Cells.Item(25, 5).Value = 12.07
$ws.Cells.Item(25, 6).Value = 10.34
$ws.Cells.Item(25, 7).Value = 13.81

$ws.Cells.Item(26, 2).Value = "Symptom – Nervous"
$ws.Cells.Item(26, 3).Value = 149
$ws.Cells.Item(26, 5).Value = 11.04
$ws.Cells.Item(26, 6).Value = 9.369999999999999
$ws.Cells.Item(26, 7).Value = 12.71

$ws.Cells.Item(27, 2).Value = "Injuries & adverse effects"
$ws.Cells.Item(27, 3).Value = 140
$ws.Cells.Item(27, 5).Value = 10.37
$ws.Cells.Item(27, 6).Value = 8.74
$ws.Cells.Item(27, 7).Value = 12

$ws.Cells.Item(28, 2).Value = "Symptom – Circulatory"
$ws.Cells.Item(28, 3).Value = 115
$ws.Cells.Item(28, 5).Value = 8.52
$ws.Cells.Item(28, 6).Value = 7.03
$ws.Cells.Item(28, 7).Value = 10.01

$ws.Cells.Item(29, 2).Value = "Other"
$ws.Cells.Item(29, 3).Value = 104
$ws.Cells.Item(29, 5).Value = 7.7
$ws.Cells.Item(29, 6).Value = 6.28
$ws.Cells.Item(29, 7).Value = 9.130000000000001

$ws.Cells.Item(30, 2).Value = "Symptom – General"
$ws.Cells.Item(30, 3).Value = 70
$ws.Cells.Item(30, 5).Value = 5.19
$ws.Cells.Item(30, 6).Value = 4
$ws.Cells.Item(30, 7).Value = 6.37

$ws.Cells.Item(31, 2).Value = "Diseases (patient-stated)"
$ws.Cells.Item(31, 3).Value = 60
$ws.Cells.Item(31, 5).Value = 4.44
$ws.Cells.Item(31, 6).Value = 3.35
$ws.Cells.Item(31, 7).Value = 5.54

$ws.Cells.Item(32, 2).Value = "Symptom – Skin/Hair/Nails"
$ws.Cells.Item(32, 3).Value = 34
$ws.Cells.Item(32, 5).Value = 2.52
$ws.Cells.Item(32, 6).Value = 1.68
$ws.Cells.Item(32, 7).Value = 3.35

$ws.Cells.Item(33, 2).Value = "Uncodable/Unknown"
$ws.Cells.Item(33, 3).Value = 29
$ws.Cells.Item(33, 5).Value = 2.15
$ws.Cells.Item(33, 6).Value = 1.37
$ws.Cells.Item(33, 7).Value = 2.92

$ws.Cells.Item(34, 2).Value = "Symptom – Musculoskeletal"
$ws.Cells.Item(34, 3).Value = 26
$ws.Cells.Item(34, 5).Value = 1.93
$ws.Cells.Item(34, 6).Value = 1.19
$ws.Cells.Item(34, 7).Value = 2.66

$ws.Cells.Item(35, 2).Value = "Symptom – Respiratory"
$ws.Cells.Item(35, 3).Value = 3354
$ws.Cells.Item(35, 5).Value = 29.66
$ws.Cells.Item(35, 6).Value = 28.82
$ws.Cells.Item(35, 7).Value = 30.5

$ws.Cells.Item(36, 2).Value = "Injuries & adverse effects"
$ws.Cells.Item(36, 3).Value = 1581
$ws.Cells.Item(36, 5).Value = 13.98
$ws.Cells.Item(36, 6).Value = 13.34
$ws.Cells.Item(36, 7).Value = 14.62

$ws.Cells.Item(37, 2).Value = "Symptom – Digestive"
$ws.Cells.Item(37, 3).Value = 1374
$ws.Cells.Item(37, 5).Value = 12.15
$ws.Cells.Item(37, 6).Value = 11.55
$ws.Cells.Item(37, 7).Value = 12.75

$ws.Cells.Item(38, 2).Value = "Symptom – Nervous"
$ws.Cells.Item(38, 3).Value = 1254
$ws.Cells.Item(38, 5).Value = 11.09
$ws.Cells.Item(38, 6).Value = 10.51
$ws.Cells.Item(38, 7).Value = 11.67

$ws.Cells.Item(39, 2).Value = "Symptom – Circulatory"
$ws.Cells.Item(39, 3).Value = 1083
$ws.Cells.Item(39, 5).Value = 9.58
$ws.Cells.Item(39, 6).Value = 9.029999999999999
$ws.Cells.Item(39, 7).Value = 10.12

$ws.Cells.Item(40, 2).Value = "Other"
$ws.Cells.Item(40, 3).Value = 742
$ws.Cells.Item(40, 5).Value = 6.56
$ws.Cells.Item(40, 6).Value = 6.1
$ws.Cells.Item(40, 7).Value = 7.02

$ws.Cells.Item(41, 2).Value = "Diseases (patient-stated)"
$ws.Cells.Item(41, 3).Value = 623
$ws.Cells.Item(41, 5).Value = 5.51
$ws.Cells.Item(41, 6).Value = 5.09
$ws.Cells.Item(41, 7).Value = 5.93

$ws.Cells.Item(42, 2).Value = "Symptom – General"
$ws.Cells.Item(42, 3).Value = 521
$ws.Cells.Item(42, 5).Value = 4.61
$ws.Cells.Item(42, 6).Value = 4.22
$ws.Cells.Item(42, 7).Value = 4.99

$ws.Cells.Item(43, 2).Value = "Symptom – Skin/Hair/Nails"
$ws.Cells.Item(43, 3).Value = 288
$ws.Cells.Item(43, 5).Value = 2.55
$ws.Cells.Item(43, 6).Value = 2.26
$ws.Cells.Item(43, 7).Value = 2.84

$ws.Cells.Item(44, 2).Value = "Symptom – Musculoskeletal"
$ws.Cells.Item(44, 3).Value = 260
$ws.Cells.Item(44, 5).Value = 2.3
$ws.Cells.Item(44, 6).Value = 2.02
$ws.Cells.Item(44, 7).Value = 2.58

$ws.Cells.Item(45, 2).Value = "Uncodable/Unknown"
$ws.Cells.Item(45, 3).Value = 229
$ws.Cells.Item(45, 5).Value = 2.02
$ws.Cells.Item(45, 6).Value = 1.77
$ws.Cells.Item(45, 7).Value = 2.28

$ws.Cells.Item(46, 2).Value = "Symptom – Respiratory"
$ws.Cells.Item(46, 3).Value = 2264
$ws.Cells.Item(46, 5).Value = 36.26
$ws.Cells.Item(46, 6).Value = 35.07
$ws.Cells.Item(46, 7).Value = 37.45

$ws.Cells.Item(47, 2).Value = "Symptom – Nervous"
$ws.Cells.Item(47, 3).Value = 760
$ws.Cells.Item(47, 5).Value = 12.17
$ws.Cells.Item(47, 6).Value = 11.36
$ws.Cells.Item(47, 7).Value = 12.98

$ws.Cells.Item(48, 2).Value = "Symptom – Digestive"
$ws.Cells.Item(48, 3).Value = 713
$ws.Cells.Item(48, 5).Value = 11.42
$ws.Cells.Item(48, 6).Value = 10.63
$ws.Cells.Item(48, 7).Value = 12.21

$ws.Cells.Item(49, 2).Value = "Injuries & adverse effects"
$ws.Cells.Item(49, 3).Value = 690
$ws.Cells.Item(49, 5).Value = 11.05
$ws.Cells.Item(49, 6).Value = 10.27
$ws.Cells.Item(49, 7).Value = 11.83

$ws.Cells.Item(50, 2).Value = "Symptom – Circulatory"
$ws.Cells.Item(50, 3).Value = 526
$ws.Cells.Item(50, 5).Value = 8.42
$ws.Cells.Item(50, 6).Value = 7.74
$ws.Cells.Item(50, 7).Value = 9.109999999999999

$ws.Cells.Item(51, 2).Value = "Other"
$ws.Cells.Item(51, 3).Value = 409
$ws.Cells.Item(51, 5).Value = 6.55
$ws.Cells.Item(51, 6).Value = 5.94
$ws.Cells.Item(51, 7).Value = 7.16

$ws.Cells.Item(52, 2).Value = "Symptom – General"
$ws.Cells.Item(52, 3).Value = 309
$ws.Cells.Item(52, 5).Value = 4.95
$ws.Cells.Item(52, 6).Value = 4.41
$ws.Cells.Item(52, 7).Value = 5.49

$ws.Cells.Item(53, 2).Value = "Diseases (patient-stated)"
$ws.Cells.Item(53, 3).Value = 255
$ws.Cells.Item(53, 5).Value = 4.08
$ws.Cells.Item(53, 6).Value = 3.59
$ws.Cells.Item(53, 7).Value = 4.57

$ws.Cells.Item(54, 2).Value = "Symptom – Skin/Hair/Nails"
$ws.Cells.Item(54, 3).Value = 118
$ws.Cells.Item(54, 5).Value = 1.89
$ws.Cells.Item(54, 6).Value = 1.55
$ws.Cells.Item(54, 7).Value = 2.23

$ws.Cells.Item(55, 2).Value = "Symptom – Musculoskeletal"
$ws.Cells.Item(55, 3).Value = 117
$ws.Cells.Item(55, 5).Value = 1.87
$ws.Cells.Item(55, 6).Value = 1.54
$ws.Cells.Item(55, 7).Value = 2.21

$ws.Cells.Item(56, 2).Value = "Uncodable/Unknown"
$ws.Cells.Item(56, 3).Value = 83
$ws.Cells.Item(56, 5).Value = 1.33
$ws.Cells.Item(56, 6).Value = 1.05
$ws.Cells.Item(56, 7).Value = 1.61
